# Update "想去人数" (interest count) values that changed between scrapes.
# Sheet "展览" (sheetId 1)
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 740
$wsExhibit.Range("F4").Value = 1473
$wsExhibit.Range("F8").Value = 6185
$wsExhibit.Range("F12").Value = 5070
$wsExhibit.Range("F13").Value = 26
$wsExhibit.Range("F14").Value = 176
$wsExhibit.Range("F22").Value = 3590
$wsExhibit.Range("F23").Value = 150

# Sheet "全部类型" (sheetId 4) - same events duplicated across the combined sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 740
$wsAll.Range("F5").Value = 1473
$wsAll.Range("F9").Value = 6185
$wsAll.Range("F13").Value = 5070
$wsAll.Range("F14").Value = 26
$wsAll.Range("F15").Value = 176
$wsAll.Range("F23").Value = 3590
$wsAll.Range("F25").Value = 150
